$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change algorithm for divide members:
# Column B and K used to hold (row - 2). Now they hold (52 - row),
# i.e. values count down from 50 to -49 instead of counting up from 0 to 99.
for ($r = 2; $r -le 101; $r++) {
    $newValue = 52 - $r
    $ws.Cells.Item($r, 2).Value = $newValue
    $ws.Cells.Item($r, 11).Value = $newValue
}

# Update the active cell/selection on the sheet view.
$ws.Range("O90").Select()
